# Weekly fruit/vegetable price update: a new week's price-record row is
# inserted at row 647 (pushing the existing rows 647-677 down to 648-678),
# and the new row is populated with that week's Espinaca (spinach) data
# for Mercado Mayorista Lo Valledor de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 647, shifting rows 647:677 down to 648:678.
$ws.Rows(647).Insert()

# Populate the newly inserted row with the new week's data.
$ws.Range("A647").Value = 6
$ws.Range("B647").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C647").Value = "Metropolitana"
$ws.Range("D647").Value = 44939
$ws.Range("E647").Value = 13
$ws.Range("F647").Value = 100112012
$ws.Range("G647").Value = "Espinaca"
$ws.Range("H647").Value = "Sin especificar"
$ws.Range("I647").Value = "Primera"
$ws.Range("J647").Value = 630
$ws.Range("K647").Value = 5000
$ws.Range("L647").Value = 6000
$ws.Range("M647").Value = 5397
$ws.Range("N647").Value = "`$/cuna 10 kilos"
$ws.Range("O647").Value = "Región Metropolitana"
$ws.Range("P647").Value = 540
$ws.Range("Q647").Value = 10
$ws.Range("R647").Value = "Hortaliza"
